$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44313
$ws.Range("M2").Value = 120

# Row 4
$ws.Range("D4").Value = 44309
$ws.Range("M4").Value = 80
$ws.Range("Q4").Value = "$/caja 14 kilos granel"
$ws.Range("S4").Value = 821
$ws.Range("T4").Value = 14

# Row 5
$ws.Range("D5").Value = 44323
$ws.Range("M5").Value = 80

# Row 7
$ws.Range("D7").Value = 44330
$ws.Range("M7").Value = 60
$ws.Range("Q7").Value = "$/caja 10 kilos empedrada"
$ws.Range("S7").Value = 11500
$ws.Range("T7").Value = 1

# Row 8
$ws.Range("D8").Value = 44322
$ws.Range("M8").Value = 60

# Row 9
$ws.Range("D9").Value = 44316
$ws.Range("M9").Value = 120

# Row 10
$ws.Range("D10").Value = 44306
